$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-08 -> 2023-09-09) for every data row (rows 2 through 439).
$ws.Range("C2:C439").Value = 45178
